$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column G (VQSR_ANNOTATIONS_SNP) is the same annotation string on every
# data row; column H (VQSR_ANNOTATIONS_INDEL) varies per row. The original
# cells stored these GATK "-an ..." argument strings wrapped in literal
# double quotes. Re-enter them without the literal quote characters, using
# a leading apostrophe so Excel applies a quote-prefix ('Text) cell format
# instead of storing the quote as part of the value.

$snpAnnotation = "-an QD -an MQ -an MQRankSum -an ReadPosRankSum -an FS -an SOR -an DP"

$indelAnnotations = @{
    2  = "-an QD -an DP -an FS -an SOR -an ReadPosRankSum -an MQRankSum"
    3  = "-an DP -an FS -an SOR -an ReadPosRankSum -an MQRankSum"
    4  = "-an QD -an FS -an SOR -an ReadPosRankSum -an MQRankSum"
    5  = "-an QD -an DP -an SOR -an ReadPosRankSum -an MQRankSum"
    6  = "-an QD -an DP -an FS -an ReadPosRankSum -an MQRankSum"
    7  = "-an QD -an DP -an FS -an SOR -an MQRankSum"
    8  = "-an QD -an DP -an FS -an SOR -an ReadPosRankSum"
    9  = "-an QD -an DP -an FS -an SOR -an ReadPosRankSum -an MQRankSum"
    10 = "-an QD -an DP -an FS -an SOR -an ReadPosRankSum -an MQRankSum"
    11 = "-an QD -an DP -an FS -an SOR -an ReadPosRankSum -an MQRankSum"
    12 = "-an QD -an DP -an FS -an SOR -an ReadPosRankSum -an MQRankSum"
    13 = "-an QD -an DP -an FS -an SOR -an ReadPosRankSum -an MQRankSum"
    14 = "-an QD -an DP -an FS -an SOR -an ReadPosRankSum -an MQRankSum"
    15 = "-an QD -an DP -an FS -an SOR -an ReadPosRankSum -an MQRankSum"
    16 = "-an QD -an DP -an FS -an SOR -an ReadPosRankSum -an MQRankSum"
    17 = "-an QD -an DP -an FS -an SOR -an ReadPosRankSum -an MQRankSum"
    18 = "-an QD -an DP -an FS -an SOR -an ReadPosRankSum -an MQRankSum"
}

for ($r = 2; $r -le 18; $r++) {
    $ws.Range("G$r").Value = "'" + $snpAnnotation
    $ws.Range("H$r").Value = "'" + $indelAnnotations[$r]
}

# Column widths were auto-fit after the edit (best-fit, per the target file).
$ws.Columns.Item(1).ColumnWidth = 23.83203125
$ws.Columns.Item(2).ColumnWidth = 13
$ws.Columns.Item(3).ColumnWidth = 10
$ws.Columns.Item(4).ColumnWidth = 15.1640625
$ws.Columns.Item(5).ColumnWidth = 23.83203125
$ws.Columns.Item(6).ColumnWidth = 15.5
$ws.Columns.Item(7).ColumnWidth = 63.6640625
$ws.Columns.Item(8).ColumnWidth = 57.1640625
$ws.Columns.Item(9).ColumnWidth = 19.83203125
$ws.Columns.Item(10).ColumnWidth = 21.33203125
$ws.Columns.Item(11).ColumnWidth = 14.83203125
$ws.Columns.Item(12).ColumnWidth = 11.83203125
$ws.Columns.Item(13).ColumnWidth = 10.33203125

# Scroll the view over towards column H and leave the selection on L1, like
# in the saved file.
$ws.Application.ActiveWindow.ScrollColumn = 8
$ws.Range("L1").Select()
